$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.851.40"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.873.77"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5348"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07181"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08159"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "1.862.20"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008562"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "26.922.12"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.404"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.313"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.732"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.734"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.618"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09147"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8192"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05001"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.943"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6067"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.628"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.202"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.27%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.644"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.939"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5128"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1492"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.988"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.627"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06130"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.20%  "
